# Insert a new data row at row 80 (this pushes the existing rows 80..163
# down to 81..164, growing the sheet's used range from A1:T163 to A1:T164),
# then populate the newly inserted row with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(80).Insert()

$ws.Cells.Item(80, 1).Value = 9
$ws.Cells.Item(80, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(80, 3).Value = "Metropolitana"
$ws.Cells.Item(80, 4).Value = 44586
$ws.Cells.Item(80, 5).Value = 13
$ws.Cells.Item(80, 6).Value = "Fruta"
$ws.Cells.Item(80, 7).Value = 100101
$ws.Cells.Item(80, 8).Value = "Berries"
$ws.Cells.Item(80, 9).Value = 100101001
$ws.Cells.Item(80, 10).Value = "Arándano (blue)"
$ws.Cells.Item(80, 11).Value = "Sin especificar"
$ws.Cells.Item(80, 12).Value = "Primera"
$ws.Cells.Item(80, 13).Value = 560
$ws.Cells.Item(80, 14).Value = 3800
$ws.Cells.Item(80, 15).Value = 4000
$ws.Cells.Item(80, 16).Value = 3921
$ws.Cells.Item(80, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(80, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(80, 19).Value = 1960
$ws.Cells.Item(80, 20).Value = 2
